# Fixing names to match
# 1) Correct the "Baseball Savant Name" column (AA) for Albert Suarez / Randy Vasquez
#    rows so they use accented spellings, consistent with the other accented
#    Baseball Savant Name entries already present in the sheet.
# 2) Multiply the raw counting-stat columns (N,O,P,T,U,V,W,X,Y,Z) for the season
#    rows (2-31) by 100 to correct a units/scale issue.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: fix Baseball Savant Name (column AA) accents ---
$ws.Cells.Item(42, 27).Value = "Albert Suárez"
$ws.Cells.Item(43, 27).Value = "Albert Suárez"

$ws.Cells.Item(48, 27).Value = "Randy Vásquez"
$ws.Cells.Item(49, 27).Value = "Randy Vásquez"

# --- Part 2: scale columns N,O,P,T,U,V,W,X,Y,Z for rows 2-31 by 100 ---
$cols = @(14, 15, 16, 20, 21, 22, 23, 24, 25, 26)

for ($row = 2; $row -le 31; $row++) {
    foreach ($col in $cols) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $cell.Value2 * 100
    }
}
